$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric must be forced to Text format first,
# otherwise Excel would auto-convert them on assignment (losing exact
# display such as trailing zeros, e.g. "239.00", or mangling the dotted
# thousands-style values). This mirrors the source data, where every cell
# in the table is stored as text.
$textCells = @(
    "D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D14", "D15",
    "D16", "D17", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27",
    "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38",
    "D40", "D41", "D44", "D45", "D48", "D49", "D50", "D51"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values (prices in column D, 1h volume % in column E,
# plus the Algorand/EnergySwap row swap in columns B/C).
$ws.Range("D2").Value = '29.320.78'
$ws.Range("D3").Value = '1.839.20'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '239.00'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").Value = '0.6219'
$ws.Range("E6").Value = '  -0.95%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '0.07326'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("D9").Value = '0.2875'
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").Value = '24.64'
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("D11").Value = '0.07729'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").Value = '1.850.29'
$ws.Range("E12").Value = '  +0.52%  '
$ws.Range("D13").Value = '4.939'
$ws.Range("E13").Value = '  -0.90%  '
$ws.Range("D14").Value = '0.6587'
$ws.Range("E14").Value = '  -2.74%  '
$ws.Range("D15").Value = '0.00001027'
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").Value = '81.33'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '6.241'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '29.317.27'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = '236.23'
$ws.Range("E19").Value = '  +3.32%  '
$ws.Range("D20").Value = '12.18'
$ws.Range("E20").Value = '  -0.97%  '
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '7.188'
$ws.Range("E22").Value = '  -3.10%  '
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '157.17'
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").Value = '8.390'
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("D26").Value = '0.1327'
$ws.Range("E26").Value = '  -1.67%  '
$ws.Range("D27").Value = '17.18'
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("D28").Value = '0.06884'
$ws.Range("E28").Value = '  +3.77%  '
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = '1.477'
$ws.Range("E30").Value = '  -0.64%  '
$ws.Range("D31").Value = '4.012'
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("D32").Value = '3.933'
$ws.Range("D33").Value = '1.153'
$ws.Range("E33").Value = '  +1.43%  '
$ws.Range("D34").Value = '1.743'
$ws.Range("E34").Value = '  -4.95%  '
$ws.Range("D35").Value = '0.6787'
$ws.Range("E35").Value = '  -2.01%  '
$ws.Range("D36").Value = '2.581'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '0.01820'
$ws.Range("E37").Value = '  -1.96%  '
$ws.Range("D38").Value = '2.779'
$ws.Range("E38").Value = '  -1.62%  '
$ws.Range("D39").Value = '1.229.86'
$ws.Range("D40").Value = '6.649'
$ws.Range("E40").Value = '  -1.62%  '
$ws.Range("D41").Value = '0.9444'
$ws.Range("E41").Value = '  +0.82%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = '1.989.65'
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").Value = '101.06'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("D45").Value = '65.06'
$ws.Range("E45").Value = '  -0.62%  '
$ws.Range("E46").Value = '  +2.35%  '
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("D48").Value = '6.859'
$ws.Range("E48").Value = '  -2.38%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.819'
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1127'
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("D51").Value = '0.3852'
$ws.Range("E51").Value = '  -1.15%  '
